# Generate Report for Archive
#
# 1. The localization status used throughout the workbook moves from
#    "Ready for handoff" to "In Translation" (stored once as a shared
#    string and referenced by every Status / zh-cn / de-de cell).
# 2. The "zh-cn"/"de-de" status columns are narrower in the new report:
#    Overview!E:F and the Status column (col C) on the "zh-cn" / "de-de"
#    sheets shrink from ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears ------------------
foreach ($ws in $wb.Worksheets) {
    $ws.UsedRange.Replace("Ready for handoff", "In Translation") | Out-Null
}

# --- 2. Narrow the status columns --------------------------------------
# ColumnWidth is expressed in (quantized) characters of the Normal style
# font; 12.5 is the value that lands on the narrower width used by the
# new report.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
